$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "Global_variable" sheet with the new values
# ---------------------------------------------------------------------------
$gv = $wb.Worksheets.Item("Global_variable")

# Figure_extension: .png -> .svg
$gv.Range("B2").Value = ".svg"

# Project folder path
$gv.Range("B3").Value = "C:\Users\amb\PycharmProjects\SI_processing_automation"

# Input folder path (value changes, note column cleared)
$gv.Range("B4").Value = "S:\Clients\T-Z\Thor Wind Farm\02_Working\Gint_databases"
$gv.Range("C4").ClearContents() | Out-Null

# Figure folder path
$gv.Range("B5").Value = "S:\Clients\T-Z\Thor Wind Farm\02_Working\CPT-data\EOS\EOS-BH-01\CPT-fig"

# Output folder path
$gv.Range("B6").Value = "S:\Clients\T-Z\Thor Wind Farm\02_Working\CPT-data\EOS\EOS-BH-01"

# Load gINT database / Load excel files - swap booleans
$gv.Range("B8").Value = $true
$gv.Range("B9").Value = $false

# database file name
$gv.Range("B13").Value = "Thor_Fugro-2022-SI.gpj"

# SCPT location (main) - value changes; other location cleared
$gv.Range("B24").Value = "EOS-BH-01"
$gv.Range("C24").ClearContents() | Out-Null

# Update the selection on this sheet (no longer the active tab)
$gv.Range("D38").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2. Insert the new "stratigraphy_color_dict" sheet right after Global_variable
# ---------------------------------------------------------------------------
$single = $wb.Worksheets.Item("Single_plot_input")
$colorDict = $wb.Worksheets.Add($null, $gv)
$colorDict.Name = "stratigraphy_color_dict"

$colorDict.Range("A1").Value = "Units"
$colorDict.Range("B1").Value = "Color"
$colorDict.Range("A2").Value = "U40-SAND"
$colorDict.Range("B2").Value = "#E0D68D"
$colorDict.Range("A3").Value = "U46-CLAY"
$colorDict.Range("B3").Value = "#228833"
$colorDict.Range("A4").Value = "U98-CLAY"
$colorDict.Range("B4").Value = "#802659"
$colorDict.Range("A5").Value = "U98-SAND"
$colorDict.Range("B5").Value = "#F1D1E3"

$colorDict.Columns.Item(1).ColumnWidth = 9.17

$colorDict.Range("B12").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. Update selections on the remaining sheets
# ---------------------------------------------------------------------------
$single.Range("F40").Select() | Out-Null

$sideBySide = $wb.Worksheets.Item("Side_by_side_plot_input")
$sideBySide.Range("M31").Select() | Out-Null

# ---------------------------------------------------------------------------
# Make the new sheet the active one (matches target workbook activeTab=1)
# ---------------------------------------------------------------------------
$colorDict.Activate() | Out-Null
$colorDict.Range("B12").Select() | Out-Null
